$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6451310861423221
$ws1.Range("C2").Value = 0.589803012746234
$ws1.Range("D2").Value = 0.9531835205992509
$ws1.Range("E2").Value = 0.7287043664996421
$ws1.Range("F2").Value = 0.8486162054018006
$ws1.Range("G2").Value = 0.9311193977344684
$ws1.Range("H2").Value = 0.8142998919889464
$ws1.Range("I2").Value = 509
$ws1.Range("J2").Value = 354
$ws1.Range("K2").Value = 180
$ws1.Range("L2").Value = 25

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.8780487804878049
$ws2.Range("C2").Value = 0.3370786516853932
$ws2.Range("D2").Value = 0.4871447902571042

$ws2.Range("B3").Value = 0.589803012746234
$ws2.Range("C3").Value = 0.9531835205992509
$ws2.Range("D3").Value = 0.7287043664996421

$ws2.Range("B4").Value = 0.6451310861423221
$ws2.Range("C4").Value = 0.6451310861423221
$ws2.Range("D4").Value = 0.6451310861423221
$ws2.Range("E4").Value = 0.6451310861423221

$ws2.Range("B5").Value = 0.7339258966170195
$ws2.Range("C5").Value = 0.6451310861423221
$ws2.Range("D5").Value = 0.6079245783783731

$ws2.Range("B6").Value = 0.7339258966170196
$ws2.Range("C6").Value = 0.6451310861423221
$ws2.Range("D6").Value = 0.6079245783783732

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 180
$ws3.Range("C2").Value = 354

$ws3.Range("B3").Value = 25
$ws3.Range("C3").Value = 509
